$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: cardholder name / account number ---
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 21.03.2025"

# --- Row 6 (was row 6: Vodafone -> now Stadtwerke Rosenheim) ---
$ws.Range("B6").Value = "23.03."
$ws.Range("C6").Value = "24.03."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 8886542"
$ws.Range("E6").Value = "87,05-"

# --- Row 7 (was Allianz -> now ZEUS BODYPOWER) ---
$ws.Range("B7").Value = "27.03."
$ws.Range("C7").Value = "28.03."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "24,66-"

# --- Row 8 (was ZEUS BODYPOWER -> now BURGER KING Backnang) ---
$ws.Range("B8").Value = "29.03."
$ws.Range("C8").Value = "30.03."
$ws.Range("D8").Value = "BURGER KING Backnang"
$ws.Range("E8").Value = "39,50-"

# --- Rows 9, 10, 11 are removed (statement now has only 3 transaction rows) ---
$ws.Range("B9:D9").ClearContents()
$ws.Range("B10:D10").ClearContents()
$ws.Range("B11:D11").ClearContents()

$ws.Range("E9").ClearContents()
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

$ws.Range("E10").ClearContents()
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

$ws.Range("E11").ClearContents()
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# --- Closing balance / next billing date ---
$ws.Range("D12").Value = "KONTOSTAND AM 02.04.2025"
$ws.Range("E12").Value = "151,21-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.04.2025"
